# ShiftSchedule_Upload.xlsx edit script
# - Updates employee ID in A2
# - Updates a handful of day-status cells in row 2 (WO <-> GN)
# - Re-points the frozen pane's top-left visible cell and active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data row (row 2) ---
$ws.Range("A2").Value = 20005312

$ws.Range("C2").Value = "WO"
$ws.Range("F2").Value = "GN"
$ws.Range("K2").Value = "WO"
$ws.Range("M2").Value = "GN"
$ws.Range("R2").Value = "WO"
$ws.Range("T2").Value = "GN"
$ws.Range("Y2").Value = "WO"
$ws.Range("AA2").Value = "GN"
$ws.Range("AF2").Value = "WO"

# --- Update the active pane selection (frozen pane's topLeftCell tracks
#     along automatically when the visible/active cell moves) ---
$ws.Range("J4").Select()
